$wb = $excel.ActiveWorkbook

# The raw width requested by the diff (40 "characters") gets a constant
# ~0.8333 padding applied by this engine's ColumnWidth setter before it is
# serialised back to the OOXML "width" attribute, so back the padding out
# here to land on an on-disk value of exactly 40.
$targetWidth = 40 - (5/6)

$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d6384f736f552b986e895a3e2459306d088e6679/e2e/26b10140-4aa8-4d1b-953b-ec8381ffae0e.md"
$hyperlinkDisplay = "26b10140-4aa8-4d1b-953b-ec8381ffae0e.md"

# Cornflower blue (FF6495ED) expressed the way the COM Color property wants
# it (BGR-packed integer), so the run's font matches the existing hyperlink
# style used elsewhere in the workbook (e.g. column A).
$hyperlinkColor = 15570276

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Latest Target File / Latest Handback File / Latest Handback DateTime /
    # Error Detail columns get a lot wider now that they carry real content.
    $ws.Columns.Item(9).ColumnWidth = $targetWidth
    $ws.Columns.Item(10).ColumnWidth = $targetWidth
    $ws.Columns.Item(16).ColumnWidth = $targetWidth

    # Row 6 is the 26b10140-... handback row; record the handback that just
    # came in even though it turned out to be stale.
    $ws.Hyperlinks.Add($ws.Range("I6"), $hyperlinkUrl, $null, $null, $hyperlinkDisplay)
    $ws.Range("I6").Font.Underline = 2
    $ws.Range("I6").Font.Color = $hyperlinkColor
}

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("J6").Value = "26b10140-4aa8-4d1b-953b-ec8381ffae0e.4c3ab74a831d3883268be8a5a4b3b94da8c84ab4.zh-cn.xlf"
$wsZh.Range("K6").Value = "2016-10-24 09:05:09"
$wsZh.Range("P6").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/98ff6fcbb926447ca6d93856fbe3f45a784d1b93/e2e/26b10140-4aa8-4d1b-953b-ec8381ffae0e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d6384f736f552b986e895a3e2459306d088e6679/e2e/26b10140-4aa8-4d1b-953b-ec8381ffae0e.md."

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("J6").Value = "26b10140-4aa8-4d1b-953b-ec8381ffae0e.4c3ab74a831d3883268be8a5a4b3b94da8c84ab4.de-de.xlf"
$wsDe.Range("K6").Value = "2016-10-24 09:05:27"
$wsDe.Range("P6").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/98ff6fcbb926447ca6d93856fbe3f45a784d1b93/e2e/26b10140-4aa8-4d1b-953b-ec8381ffae0e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d6384f736f552b986e895a3e2459306d088e6679/e2e/26b10140-4aa8-4d1b-953b-ec8381ffae0e.md."
